# Auto-generated edit script applying the scheduled-runner profit/price updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit worksheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value2 = 13500
$ws.Range("J3").Value2 = 13500
$ws.Range("L3").Value2 = 13500
$ws.Range("N3").Value2 = -13728
# Row 40
$ws.Range("H40").Value2 = 1300
$ws.Range("I40").Value2 = 1100
$ws.Range("K40").Value2 = 1100
$ws.Range("M40").Value2 = -925
# Row 98
$ws.Range("H98").Value2 = 3247.0527
$ws.Range("I98").Value2 = 3930.3076
$ws.Range("J98").Value2 = 1766.6666
$ws.Range("K98").Value2 = 3930.3076
$ws.Range("L98").Value2 = 1766.6666
$ws.Range("M98").Value2 = -2432.3076
$ws.Range("N98").Value2 = -4762.6666
# Row 102
$ws.Range("H102").Value2 = 13500
$ws.Range("J102").Value2 = 13500
$ws.Range("L102").Value2 = 13500
$ws.Range("N102").Value2 = -19990
# Row 122
$ws.Range("H122").Value2 = 3247.0527
$ws.Range("I122").Value2 = 3930.3076
$ws.Range("J122").Value2 = 1766.6666
$ws.Range("K122").Value2 = 11790.9228
$ws.Range("L122").Value2 = 5299.9998
$ws.Range("M122").Value2 = -9340.9228
$ws.Range("N122").Value2 = -10199.9998
# Row 132
$ws.Range("H132").Value2 = 388393.22
$ws.Range("I132").Value2 = 438835.8
$ws.Range("J132").Value2 = 1666.6666
$ws.Range("K132").Value2 = 1316507.4
$ws.Range("L132").Value2 = 4999.9998
$ws.Range("M132").Value2 = -1313977.4
$ws.Range("N132").Value2 = -10059.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 8265.617
$ws.Range("I32").Value2 = 9536.556
$ws.Range("K32").Value2 = 9536.556
$ws.Range("M32").Value2 = -9249.556
# Row 61
$ws.Range("H61").Value2 = 4257.6924
$ws.Range("I61").Value2 = 3887
$ws.Range("J61").Value2 = 4575.4287
$ws.Range("K61").Value2 = 3887
$ws.Range("L61").Value2 = 4575.4287
$ws.Range("M61").Value2 = -3675
$ws.Range("N61").Value2 = -4999.4287
# Row 62
$ws.Range("H62").Value2 = 0
$ws.Range("J62").Value2 = 0
$ws.Range("L62").Value2 = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("N65").ClearContents()
# Row 136
$ws.Range("H136").Value2 = 4257.6924
$ws.Range("I136").Value2 = 3887
$ws.Range("J136").Value2 = 4575.4287
$ws.Range("K136").Value2 = 11661
$ws.Range("L136").Value2 = 13726.2861
$ws.Range("M136").Value2 = -9111
$ws.Range("N136").Value2 = -18826.2861

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 76
$ws.Range("H76").Value2 = 23157
$ws.Range("I76").Value2 = 12000
$ws.Range("J76").Value2 = 34314
$ws.Range("K76").Value2 = 12000
$ws.Range("L76").Value2 = 34314
$ws.Range("M76").Value2 = -11685
$ws.Range("N76").Value2 = -34944
# Row 79
$ws.Range("H79").Value2 = 23157
$ws.Range("I79").Value2 = 12000
$ws.Range("J79").Value2 = 34314
$ws.Range("K79").Value2 = 12000
$ws.Range("L79").Value2 = 34314
$ws.Range("M79").Value2 = -10908
$ws.Range("N79").Value2 = -36498
# Row 86
$ws.Range("H86").Value2 = 2027.1724
$ws.Range("I86").Value2 = 1773.3914
$ws.Range("J86").Value2 = 3000
$ws.Range("K86").Value2 = 1773.3914
$ws.Range("L86").Value2 = 3000
$ws.Range("M86").Value2 = -650.3914
$ws.Range("N86").Value2 = -5246
# Row 89
$ws.Range("H89").Value2 = 2027.1724
$ws.Range("I89").Value2 = 1773.3914
$ws.Range("J89").Value2 = 3000
$ws.Range("K89").Value2 = 8866.957
$ws.Range("L89").Value2 = 15000
$ws.Range("M89").Value2 = -3250.957
$ws.Range("N89").Value2 = -26232
# Row 94
$ws.Range("H94").Value2 = 1041.8889
$ws.Range("I94").Value2 = 1053.72
$ws.Range("K94").Value2 = 1053.72
$ws.Range("M94").Value2 = -602.72
# Row 99
$ws.Range("H99").Value2 = 1549.8572
$ws.Range("J99").Value2 = 1933.3334
$ws.Range("L99").Value2 = 1933.3334
$ws.Range("N99").Value2 = -4929.3334
# Row 107
$ws.Range("H107").Value2 = 1498.8148
$ws.Range("I107").Value2 = 1429.7106
$ws.Range("K107").Value2 = 1429.7106
$ws.Range("M107").Value2 = 490.2893999999999
# Row 134
$ws.Range("H134").Value2 = 59714.42
$ws.Range("I134").Value2 = 66489.35000000001
$ws.Range("J134").Value2 = 2127.5
$ws.Range("K134").Value2 = 199468.05
$ws.Range("L134").Value2 = 6382.5
$ws.Range("M134").Value2 = -196933.05
$ws.Range("N134").Value2 = -11452.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value2 = 1192.3334
$ws.Range("I99").Value2 = 1281.5
$ws.Range("J99").Value2 = 1014
$ws.Range("K99").Value2 = 1281.5
$ws.Range("L99").Value2 = 1014
$ws.Range("M99").Value2 = 216.5
$ws.Range("N99").Value2 = -4010
# Row 126
$ws.Range("H126").Value2 = 1192.3334
$ws.Range("I126").Value2 = 1281.5
$ws.Range("J126").Value2 = 1014
$ws.Range("K126").Value2 = 3844.5
$ws.Range("L126").Value2 = 3042
$ws.Range("M126").Value2 = -1374.5
$ws.Range("N126").Value2 = -7982

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value2 = 601
$ws.Range("I113").Value2 = 621.2857
$ws.Range("J113").Value2 = 591.5333000000001
$ws.Range("K113").Value2 = 1863.8571
$ws.Range("L113").Value2 = 1774.5999
$ws.Range("M113").Value2 = 306.1428999999998
$ws.Range("N113").Value2 = -6114.5999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value2 = 431177.84
$ws.Range("J21").Value2 = 2001330
$ws.Range("L21").Value2 = 2001330
$ws.Range("N21").Value2 = -2001676
# Row 30
$ws.Range("H30").Value2 = 431177.84
$ws.Range("J30").Value2 = 2001330
$ws.Range("L30").Value2 = 2001330
$ws.Range("N30").Value2 = -2001540
# Row 39
$ws.Range("H39").Value2 = 17252.2
$ws.Range("J39").Value2 = 17252.2
$ws.Range("L39").Value2 = 17252.2
$ws.Range("N39").Value2 = -18316.2
# Row 80
$ws.Range("H80").Value2 = 135255.11
$ws.Range("J80").Value2 = 151787
$ws.Range("L80").Value2 = 151787
$ws.Range("N80").Value2 = -153783
# Row 83
$ws.Range("H83").Value2 = 135255.11
$ws.Range("J83").Value2 = 151787
$ws.Range("L83").Value2 = 758935
$ws.Range("N83").Value2 = -768919
# Row 113
$ws.Range("H113").Value2 = 2032.3636
$ws.Range("I113").Value2 = 2454.3333
$ws.Range("J113").Value2 = 1526
$ws.Range("K113").Value2 = 2454.3333
$ws.Range("L113").Value2 = 1526
$ws.Range("M113").Value2 = -284.3332999999998
$ws.Range("N113").Value2 = -5866

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value2 = 1028.2941
$ws.Range("I46").Value2 = 1111.5834
$ws.Range("K46").Value2 = 1111.5834
$ws.Range("M46").Value2 = -923.5834
# Row 61
$ws.Range("H61").Value2 = 1142.4286
$ws.Range("J61").Value2 = 1333.3334
$ws.Range("L61").Value2 = 1333.3334
$ws.Range("N61").Value2 = -1737.3334
# Row 93
$ws.Range("H93").Value2 = 1633.3334
$ws.Range("I93").Value2 = 1525
$ws.Range("J93").Value2 = 2500
$ws.Range("K93").Value2 = 1525
$ws.Range("L93").Value2 = 2500
$ws.Range("M93").Value2 = -277
$ws.Range("N93").Value2 = -4996
# Row 106
$ws.Range("H106").Value2 = 19000
$ws.Range("J106").Value2 = 19000
$ws.Range("L106").Value2 = 19000
$ws.Range("N106").Value2 = -21524
# Row 113
$ws.Range("H113").Value2 = 1142.4286
$ws.Range("J113").Value2 = 1333.3334
$ws.Range("L113").Value2 = 1333.3334
$ws.Range("N113").Value2 = -5673.3334
# Row 122
$ws.Range("H122").Value2 = 5595.324
$ws.Range("I122").Value2 = 6374.4287
$ws.Range("J122").Value2 = 3171.4443
$ws.Range("K122").Value2 = 19123.2861
$ws.Range("L122").Value2 = 9514.332900000001
$ws.Range("M122").Value2 = -16673.2861
$ws.Range("N122").Value2 = -14414.3329

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("N92").ClearContents()
# Row 98
$ws.Range("H98").Value2 = 3500
$ws.Range("J98").Value2 = 3500
$ws.Range("L98").Value2 = 3500
$ws.Range("N98").Value2 = -9490
# Row 104
$ws.Range("H104").Value2 = 22228.166
$ws.Range("J104").Value2 = 22228.166
$ws.Range("L104").Value2 = 22228.166
$ws.Range("N104").Value2 = -29216.166
# Row 107
$ws.Range("H107").Value2 = 474.46155
$ws.Range("I107").Value2 = 307.26315
$ws.Range("J107").Value2 = 928.2857
$ws.Range("K107").Value2 = 921.78945
$ws.Range("L107").Value2 = 2784.8571
$ws.Range("M107").Value2 = 998.21055
$ws.Range("N107").Value2 = -6624.8571
# Row 122
$ws.Range("H122").Value2 = 1581.3334
$ws.Range("I122").Value2 = 1332.0714
$ws.Range("J122").Value2 = 1930.3
$ws.Range("K122").Value2 = 3996.2142
$ws.Range("L122").Value2 = 5790.9
$ws.Range("M122").Value2 = -1546.2142
$ws.Range("N122").Value2 = -10690.9
# Row 126
$ws.Range("H126").Value2 = 2426.5
$ws.Range("I126").Value2 = 2609.913
$ws.Range("J126").Value2 = 1582.8
$ws.Range("K126").Value2 = 7829.739
$ws.Range("L126").Value2 = 4748.4
$ws.Range("M126").Value2 = -5359.739
$ws.Range("N126").Value2 = -9688.4

